$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric (e.g. "307.08") must be forced to
# text so they keep the original text representation (matches source data,
# which stores prices as literal strings, some with thousands separators).
$textCells = @("D5","D6","D7","D9","D10","D11","D12","D16","D17","D19","D21","D22","D23","D24","D25","D26","D27","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D44","D47","D48","D49","D50","D51")
$textRange = $ws.Range($textCells[0])
foreach ($addr in $textCells[1..($textCells.Length-1)]) {
    $textRange = $excel.Union($textRange, $ws.Range($addr))
}
foreach ($area in $textRange.Areas) {
    $area.NumberFormat = "@"
}

$ws.Range('D2').Value = '46.708.92'
$ws.Range('E2').Value = '  +4.72%  '
$ws.Range('D3').Value = '2.345.05'
$ws.Range('E3').Value = '  +4.35%  '
$ws.Range('E4').Value = '  -0.70%  '
$ws.Range('D5').Value = '307.08'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').Value = '99.33'
$ws.Range('E6').Value = '  +5.01%  '
$ws.Range('D7').Value = '0.579'
$ws.Range('E7').Value = '  +1.44%  '
$ws.Range('E8').Value = '  -0.51%  '
$ws.Range('D9').Value = '0.536'
$ws.Range('E9').Value = '  +3.65%  '
$ws.Range('D10').Value = '36.05'
$ws.Range('E10').Value = '  +3.93%  '
$ws.Range('D11').Value = '0.0808'
$ws.Range('E11').Value = '  +0.80%  '
$ws.Range('D12').Value = '7.43'
$ws.Range('E12').Value = '  +3.30%  '
$ws.Range('E13').Value = '  -0.15%  '
$ws.Range('D14').Value = '2.704.30'
$ws.Range('E14').Value = '  +4.48%  '
$ws.Range('D15').Value = '2.344.90'
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('D16').Value = '14.28'
$ws.Range('E16').Value = '  +5.36%  '
$ws.Range('D17').Value = '0.830'
$ws.Range('E17').Value = '  -0.33%  '
$ws.Range('D18').Value = '46.712.20'
$ws.Range('E18').Value = '  +5.35%  '
$ws.Range('D19').Value = '13.38'
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('D21').Value = '6.19'
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('D22').Value = '66.73'
$ws.Range('E22').Value = '  +2.34%  '
$ws.Range('D23').Value = '245.36'
$ws.Range('E23').Value = '  +3.43%  '
$ws.Range('D24').Value = '2.97'
$ws.Range('E24').Value = '  +0.73%  '
$ws.Range('D25').Value = '1.98'
$ws.Range('E25').Value = '  +0.66%  '
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.36%  '
$ws.Range('D27').Value = '41.96'
$ws.Range('E27').Value = '  +13.31%  '
$ws.Range('E28').Value = '  -0.90%  '
$ws.Range('D29').Value = '9.89'
$ws.Range('E29').Value = '  +1.37%  '
$ws.Range('D30').Value = '20.24'
$ws.Range('E30').Value = '  +1.57%  '
$ws.Range('D31').Value = '5.74'
$ws.Range('E31').Value = '  -2.49%  '
$ws.Range('D32').Value = '151.37'
$ws.Range('E32').Value = '  +1.86%  '
$ws.Range('D33').Value = '0.0814'
$ws.Range('E33').Value = '  +3.93%  '
$ws.Range('D34').Value = '2.61'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').Value = '3.08'
$ws.Range('E35').Value = '  -4.12%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = '0.108'
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('B37').Value = 'Stellar'
$ws.Range('C37').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D37').Value = '0.119'
$ws.Range('E37').Value = '  +0.49%  '
$ws.Range('D38').Value = '1.83'
$ws.Range('E38').Value = '  -1.72%  '
$ws.Range('D39').Value = '4.04'
$ws.Range('E39').Value = '  +6.89%  '
$ws.Range('D40').Value = '0.0320'
$ws.Range('E40').Value = '  +7.17%  '
$ws.Range('D41').Value = '3.44'
$ws.Range('E41').Value = '  +2.53%  '
$ws.Range('D42').Value = '13.87'
$ws.Range('E42').Value = '  -9.11%  '
$ws.Range('E43').Value = '  -0.60%  '
$ws.Range('D44').Value = '1.94'
$ws.Range('E44').Value = '  +10.57%  '
$ws.Range('D45').Value = '1.808.02'
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('E46').Value = '  +5.91%  '
$ws.Range('D47').Value = '81.36'
$ws.Range('E47').Value = '  -0.62%  '
$ws.Range('D48').Value = '73.39'
$ws.Range('E48').Value = '  +7.22%  '
$ws.Range('D49').Value = '4.92'
$ws.Range('E49').Value = '  +1.78%  '
$ws.Range('D50').Value = '98.61'
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('D51').Value = '55.59'
$ws.Range('E51').Value = '  +3.40%  '
